$d = $word.ActiveDocument

# Locate the "ACT" Heading 2 paragraph, then delete the very next paragraph
# (the italic "Acts" sub-title paragraph that immediately follows it).
$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count - 1; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Style.NameLocal -eq "Heading 2" -and $p.Range.Text.Trim() -eq "ACT") {
        $next = $d.Paragraphs.Item($i + 1)
        if ($next.Range.Text.Trim() -eq "Acts") {
            $target = $next
        }
        break
    }
}

if ($target -ne $null) {
    $target.Range.Delete()
}
